# Generate Report for Handoff
# Updates Status from "In Translation" to "Ready for handoff" and refreshes
# the "Latest Handoff Date(time)" values on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: B2 = zh-cn status, C2 = de-de status, D2 = latest handoff date
$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"
$wsOverview.Range("D2").Value = "2016-03-24 00:38:16"

# zh-cn sheet: C2 = status, E2 = latest handoff datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "2016-03-24 00:38:12"

# de-de sheet: C2 = status, E2 = latest handoff datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "2016-03-24 00:38:16"
